$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.587229
$ws.Range("H2").Value = 31.761687
$ws.Range("I2").Value = 0.1340590927938227
$ws.Range("J2").Value = 0.1340590927938227
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03927866666666666
$ws.Range("N2").Value = 0.117836
$ws.Range("O2").Value = 0.7432525340448212
$ws.Range("P2").Value = 0.7432525340448213
$ws.Range("Q2").Value = 0.4158522388146667
$ws.Range("R2").Value = 3.742670149332
$ws.Range("S2").Value = 0.09963976043075852
$ws.Range("T2").Value = 0.09963976043075853

$ws.Range("G3").Value = 10.587229
$ws.Range("H3").Value = 31.761687
$ws.Range("I3").Value = 0.1340590927938227
$ws.Range("J3").Value = 0.1340590927938227
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01356833333333333
$ws.Range("N3").Value = 0.040705
$ws.Range("O3").Value = 0.2567474659551788
$ws.Range("P3").Value = 0.2567474659551788
$ws.Range("Q3").Value = 0.1436510521483333
$ws.Range("R3").Value = 1.292859469335
$ws.Range("S3").Value = 0.03441933236306414
$ws.Range("T3").Value = 0.03441933236306414

$ws.Range("I4").Value = 0.2241137347582675
$ws.Range("J4").Value = 0.2241137347582675
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03927866666666666
$ws.Range("N4").Value = 0.117836
$ws.Range("O4").Value = 0.7432525340448212
$ws.Range("P4").Value = 0.7432525340448213
$ws.Range("Q4").Value = 0.6952023649128889
$ws.Range("R4").Value = 6.256821284216
$ws.Range("S4").Value = 0.1665731012733313
$ws.Range("T4").Value = 0.1665731012733313

$ws.Range("I5").Value = 0.2241137347582675
$ws.Range("J5").Value = 0.2241137347582675
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01356833333333333
$ws.Range("N5").Value = 0.040705
$ws.Range("O5").Value = 0.2567474659551788
$ws.Range("P5").Value = 0.2567474659551788
$ws.Range("Q5").Value = 0.2401491247477778
$ws.Range("R5").Value = 2.16134212273
$ws.Range("S5").Value = 0.05754063348493627
$ws.Range("T5").Value = 0.05754063348493627

$ws.Range("G6").Value = 34.22308866666667
$ws.Range("H6").Value = 102.669266
$ws.Range("I6").Value = 0.43334438305395
$ws.Range("J6").Value = 0.43334438305395
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03927866666666666
$ws.Range("N6").Value = 0.117836
$ws.Range("O6").Value = 0.7432525340448212
$ws.Range("P6").Value = 0.7432525340448213
$ws.Range("Q6").Value = 1.344237292041778
$ws.Range("R6").Value = 12.098135628376
$ws.Range("S6").Value = 0.322084310818938
$ws.Range("T6").Value = 0.322084310818938

$ws.Range("G7").Value = 34.22308866666667
$ws.Range("H7").Value = 102.669266
$ws.Range("I7").Value = 0.43334438305395
$ws.Range("J7").Value = 0.43334438305395
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01356833333333333
$ws.Range("N7").Value = 0.040705
$ws.Range("O7").Value = 0.2567474659551788
$ws.Range("P7").Value = 0.2567474659551788
$ws.Range("Q7").Value = 0.4643502747255556
$ws.Range("R7").Value = 4.17915247253
$ws.Range("S7").Value = 0.111260072235012
$ws.Range("T7").Value = 0.111260072235012

$ws.Range("G8").Value = 16.46479166666667
$ws.Range("H8").Value = 49.394375
$ws.Range("I8").Value = 0.2084827893939599
$ws.Range("J8").Value = 0.2084827893939599
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.03927866666666666
$ws.Range("N8").Value = 0.117836
$ws.Range("O8").Value = 0.7432525340448212
$ws.Range("P8").Value = 0.7432525340448213
$ws.Range("Q8").Value = 0.6467150636111111
$ws.Range("R8").Value = 5.820435572499999
$ws.Range("S8").Value = 0.1549553615217935
$ws.Range("T8").Value = 0.1549553615217935

$ws.Range("G9").Value = 16.46479166666667
$ws.Range("H9").Value = 49.394375
$ws.Range("I9").Value = 0.2084827893939599
$ws.Range("J9").Value = 0.2084827893939599
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01356833333333333
$ws.Range("N9").Value = 0.040705
$ws.Range("O9").Value = 0.2567474659551788
$ws.Range("P9").Value = 0.2567474659551788
$ws.Range("Q9").Value = 0.2233997815972222
$ws.Range("R9").Value = 2.010598034375
$ws.Range("S9").Value = 0.05352742787216642
$ws.Range("T9").Value = 0.05352742787216642

